$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "several new people added" -> the underlying data point in A2 grew
# from 420 to 1000. B2 (=$B$1/$A$1*A2) recalculates automatically.
$ws.Range("A2").Value = 1000

# Reflect the author's resized/repositioned workbook window, best effort.
$win = $excel.ActiveWindow
$win.Left = 13840
$win.Width = 13040
